$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.531.69'
$ws.Range('D3').Value = '1.587.40'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('E4').Value = '  +1.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.76'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.492'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('E7').Value = '  +0.98%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.27'
$ws.Range('E8').Value = '  +6.04%  '
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('E11').Value = '  +1.82%  '
$ws.Range('D12').Value = '1.814.33'
$ws.Range('E12').Value = '  +0.93%  '
$ws.Range('D13').Value = '1.600.73'
$ws.Range('E13').Value = '  +2.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.528'
$ws.Range('E14').Value = '  +1.75%  '
$ws.Range('D16').Value = '28.546.75'
$ws.Range('E16').Value = '  +3.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.05'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '232.00'
$ws.Range('E18').Value = '  +2.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.49'
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('D20').Value = '0.0₃0707'
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.04'
$ws.Range('E22').Value = '  -1.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.31'
$ws.Range('E23').Value = '  -1.11%  '
$ws.Range('E24').Value = '  +1.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.73'
$ws.Range('E25').Value = '  +0.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.23'
$ws.Range('E26').Value = '  +0.36%  '
$ws.Range('E27').Value = '  -0.93%  '
$ws.Range('E28').Value = '  -0.71%  '
$ws.Range('E30').Value = '  -1.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0469'
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.26'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('E33').Value = '  +1.16%  '
$ws.Range('D34').Value = '1.387.33'
$ws.Range('E34').Value = '  -4.79%  '
$ws.Range('E35').Value = '  -1.34%  '
$ws.Range('E36').Value = '  -10.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.62'
$ws.Range('E38').Value = '  +10.53%  '
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.540'
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('E42').Value = '  +0.96%  '
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('E45').Value = '  +0.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '62.88'
$ws.Range('E46').Value = '  -1.95%  '
$ws.Range('D47').Value = '1.724.94'
$ws.Range('E47').Value = '  +0.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.13'
$ws.Range('E48').Value = '  +1.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '87.08'
$ws.Range('E49').Value = '  +0.10%  '
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0521'
$ws.Range('E51').Value = '  -1.36%  '
